$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R (year 2021) mirrors the formatting of column Q for every row.
# Copy each Q-column cell onto the corresponding R-column cell first so the
# style (number format, borders, etc.) carries over, then overwrite with the
# 2021 values from the source data.

$values = @{
    3  = 2021
    4  = 1.7931687443515183
    5  = 1.0977143806517458
    6  = 2.4989281705678046
    7  = 2.3489023398681002
    8  = 1.8410239038543676
    9  = 2.8382683724659588
    10 = 1.2584206034913306
    11 = 0.79202525610136665
    12 = 1.7183687369364922
    13 = 1.7860084101151579
    14 = 1.5807090270340762
    15 = 1.9930959157478496
    16 = 1.0231016349164126
    17 = 0
    18 = 2.0091214112068791
    19 = 2.2092990108041848
    20 = 0.86496336159360854
    21 = 3.5236628052020538
    22 = 1.4678252700798498
    23 = 0.74155920237892192
    24 = 2.1792664589099311
    25 = 1.5302890103825006
    26 = 0.80351618683358383
    27 = 2.280288974802807
    28 = 2.3014726663297309
    29 = 1.7358308467556451
    30 = 2.9402079315049163
    31 = 1.2198989923634325
    32 = 1.1878318505232399
    33 = 1.2537455648750642
}

foreach ($row in 3..33) {
    $ws.Range("Q$row").Copy($ws.Range("R$row"))
    $ws.Range("R$row").Value = $values[$row]
}

# Move the active selection to S14, matching the post-edit workbook state.
$ws.Range("S14").Select()
